$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to the municipality name
$ws.Name = "დუშეთი"

# Remove the outdated "(census results)" caption from A2 and its companion blank cell in B2
$ws.Range("A2:B2").Clear() | Out-Null

# Remove the now-superfluous blank spacer row (old row 3)
$ws.Rows("3").Delete() | Out-Null

# Remove the stray formatted-but-empty cell next to the title
$ws.Range("B1").Clear() | Out-Null

# Drop the 1989 and 2002 columns, keeping only the latest (2014) figures
$ws.Range("C:D").Delete() | Out-Null

# Update the remaining year column to read 2014
$ws.Range("B4").Value = 2014

# Leave the selection on A2, like in the saved file
$ws.Range("A2").Select() | Out-Null
